$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B23").Value = "Atom autocomplete:"
$f = $ws.Range("B23").Font
$f.Italic = $false
$f.FontStyle = "Regular"
